$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.567.76'
$ws.Range("E2").Value = '  +2.13%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.383.83'
$ws.Range("E3").Value = '  +1.90%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '553.22'
$ws.Range("E5").Value = '  +2.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.39'
$ws.Range("E6").Value = '  +2.95%  '

$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.526'
$ws.Range("E8").Value = '  +2.34%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.386.13'
$ws.Range("E9").Value = '  +2.07%  '

$ws.Range("E10").Value = '  +6.03%  '

$ws.Range("E11").Value = '  +2.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.36'
$ws.Range("E12").Value = '  +3.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.351'
$ws.Range("E13").Value = '  +4.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.52'
$ws.Range("E14").Value = '  +3.88%  '

$ws.Range("E15").Value = '  +5.80%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '61.509.74'
$ws.Range("E16").Value = '  +2.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '11.01'
$ws.Range("E17").Value = '  +5.21%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '321.42'
$ws.Range("E18").Value = '  +2.89%  '

$ws.Range("E19").Value = '  +2.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.78'
$ws.Range("E20").Value = '  +4.81%  '

$ws.Range("E21").Value = '  -0.23%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.76'
$ws.Range("E22").Value = '  -5.25%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.50'
$ws.Range("E23").Value = '  +2.57%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.95'
$ws.Range("E24").Value = '  +11.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.25'
$ws.Range("E25").Value = '  +5.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '522.78'
$ws.Range("E26").Value = '  +4.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0₃0909'
$ws.Range("E27").Value = '  +2.78%  '

$ws.Range("E28").Value = '  +5.39%  '

$ws.Range("E29").Value = '  +2.01%  '

$ws.Range("E30").Value = '  +2.84%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.55'
$ws.Range("E31").Value = '  +2.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.17%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.59'
$ws.Range("E33").Value = '  +7.67%  '

$ws.Range("E34").Value = '  +5.94%  '

$ws.Range("E35").Value = '  +8.85%  '

$ws.Range("E36").Value = '  +3.07%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.58'
$ws.Range("E37").Value = '  +2.56%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '146.91'
$ws.Range("E38").Value = '  +5.62%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '41.40'
$ws.Range("E40").Value = '  +3.61%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '148.76'
$ws.Range("E41").Value = '  +9.16%  '

$ws.Range("E42").Value = '  +6.10%  '

$ws.Range("E43").Value = '  +3.09%  '

$ws.Range("E44").Value = '  +3.91%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.80'
$ws.Range("E45").Value = '  +2.66%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.584'
$ws.Range("E46").Value = '  +3.68%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0909'
$ws.Range("E47").Value = '  +2.12%  '

$ws.Range("E48").Value = '  +2.41%  '

$ws.Range("E49").Value = '  +0.76%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.78'
$ws.Range("E50").Value = '  +2.49%  '

$ws.Range("E51").Value = '  +5.06%  '
